# Toggle the "Started" (column C) Yes/No flag for a handful of rows on the
# "by Coach" sheet, and move the frozen pane's top-left visible row up from
# row 67 to row 63.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by Coach")

# Rows whose Column C value flips between "Yes" and "No".
$rowsToToggle = @(14, 15, 18, 21, 50, 55, 65, 67, 77, 83)

foreach ($r in $rowsToToggle) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "Yes") {
        $cell.Value = "No"
    } elseif ($cell.Value2 -eq "No") {
        $cell.Value = "Yes"
    }
}

# Scroll the frozen pane up a bit: topLeftCell moves from A67 to A63.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 63
